# Update column G ("K") values for rows 2-19 on the active worksheet,
# replacing the old Strike# counts with the newly regenerated K values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 3
    3  = 0
    4  = 2
    5  = 2
    6  = 1
    7  = 2
    8  = 5
    9  = 2
    10 = 2
    11 = 5
    12 = 4
    13 = 5
    14 = 6
    15 = 6
    16 = 4
    17 = 3
    18 = 2
    19 = 4
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
